# Add the new "ODI Bowling Extra" worksheet right after "ODI Batting Extra",
# matching the workbook.xml diff (sheetId=5, positioned as the 5th tab).
$wb = $excel.ActiveWorkbook
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $battingExtra)
$newSheet.Name = "ODI Bowling Extra"

# Reuse the exact header formatting (bold font, thin border, centered/top
# aligned) already used by the other sheets' header rows, by copying the
# format from the existing "ODI Batting Extra" header cells.
$battingExtra.Range("A1:C1").Copy() | Out-Null
$newSheet.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Header row.
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "MAIDEN_OVERS"
$newSheet.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

# Data rows. MATCH_CODE, MAIDEN_OVERS and PERCENT_WICKETS_OF_ALL are all
# stored as plain text in the source data (match codes / counts that look
# numeric, and percentages formatted as text strings), so every value is
# written with a leading apostrophe to force text storage instead of
# Excel's automatic number/percentage inference. $null entries are left
# untouched so the cell stays blank, matching the sparse source rows.
$data = @(
    @("4248", "1", "20.00%"),
    @("4249", "0", "10.00%"),
    @("4268", "0", "10.00%"),
    @("4270", "0", "30.00%"),
    @("4310", $null, $null),
    @("4316", "0", "30.00%"),
    @("4324", "0", $null),
    @("4345", "0", "10.00%"),
    @("4350", "0", "10.00%"),
    @("4353", $null, $null),
    @("4359", $null, $null),
    @("4360", "0", "40.00%"),
    @("4362", $null, $null),
    @("4454", "0", "20.00%"),
    @("4456", $null, $null),
    @("4457", "0", "30.00%"),
    @("4480", "0", $null),
    @("4482", "0", "30.00%"),
    @("4524", $null, $null),
    @("4526", "0", $null)
)

$row = 2
foreach ($entry in $data) {
    $newSheet.Cells.Item($row, 1).Value = "'" + $entry[0]
    if ($entry[1] -ne $null) {
        $newSheet.Cells.Item($row, 2).Value = "'" + $entry[1]
    }
    if ($entry[2] -ne $null) {
        $newSheet.Cells.Item($row, 3).Value = "'" + $entry[2]
    }
    $row = $row + 1
}

Write-Output ("Worksheets: " + ($wb.Worksheets | ForEach-Object { $_.Name }))

# Restore the original active sheet/selection (adding a sheet makes it the
# active tab as a side effect; the source workbook keeps "Player Info" / A1
# active).
$wb.Worksheets.Item("Player Info").Activate()
$wb.Worksheets.Item("Player Info").Range("A1").Select() | Out-Null
